# feat: add win screen transition
# Also debugged missing references in option menu tab buttons
#
# Appends two new rows (44 "Game won" / 45 "Death") to the UI defaults
# table on tsv_UI_Defaults, carrying the localized win/lose screen copy
# (English / French / Spanish / Japanese / Simplified Chinese).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 ("Dodge delay") carries the standard data-row formatting
# (bordered, vertically centered, wrapped, 28.8pt tall). Clone that
# formatting onto the two new rows before filling in values.
$ws.Range("A43:G43").Copy()
$ws.Range("A44:G45").PasteSpecial(-4122)

$ws.Rows.Item(44).RowHeight = 28.8
$ws.Rows.Item(45).RowHeight = 28.8

# Column order: A=Notes, B=Key | Language, C=English, D=Français,
# E=Español, F=日本語, G=简体中文
$ws.Range("A44").Value = "Game won"
$ws.Range("A45").Value = "Death"
$ws.Range("B44").Value = "tmp game won"
$ws.Range("B45").Value = "tmp game lost"
$ws.Range("C45").Value = "You died!"
$ws.Range("C44").Value = "You survived!"
$ws.Range("D44").Value = "Vous avez survécu!"
$ws.Range("D45").Value = "Vous êtes mort!"
$ws.Range("E45").Value = "¡Moriste!"
$ws.Range("E44").Value = "¡Sobreviviste!"
$ws.Range("F45").Value = "死にました！"
$ws.Range("F44").Value = "生き残った！"
$ws.Range("G44").Value = "你活下来了！"
$ws.Range("G45").Value = "你死了！"

$ws.Range("F47").Select() | Out-Null
